$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 151, pushing existing rows 151:228 down to 152:229.
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with the new Ajo price record.
$ws.Range("A151").Value = 5
$ws.Range("B151").Value = "Macroferia Regional de Talca"
$ws.Range("C151").Value = "Maule"
$ws.Range("D151").Value = 44572
$ws.Range("E151").Value = 7
$ws.Range("F151").Value = 100112003
$ws.Range("G151").Value = "Ajo"
$ws.Range("H151").Value = "Chino"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 300
$ws.Range("K151").Value = 21000
$ws.Range("L151").Value = 21000
$ws.Range("M151").Value = 21000
$ws.Range("N151").Value = "$/malla 10 kilos"
$ws.Range("O151").Value = "China"
$ws.Range("P151").Value = 2100
$ws.Range("Q151").Value = 10
$ws.Range("R151").Value = "Hortaliza"
